# Add "godot" as a new repository row (row 10) to both distribution sheets.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "all branch segments"
$ws2 = $wb.Worksheets.Item(2)   # "only branch segs gt 1"

# Sheet 1: "all branch segments" - add row 10 for godot
$ws1.Range("A10").Value = "godot"
$ws1.Range("B10").Value = 1.5248756218999999
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 1.2596275068
$ws1.Range("E10").Value = 26.654228856
$ws1.Range("F10").Value = 2
$ws1.Range("G10").Value = 184.99715531999999
$ws1.Range("H10").Value = 14.650171759999999
$ws1.Range("I10").Value = 1
$ws1.Range("J10").Value = 82.991788528000001
$ws1.Range("K10").Value = 26.661691542
$ws1.Range("L10").Value = 2
$ws1.Range("M10").Value = 184.99639404999999
$ws1.Range("N10").Value = 14.653156835000001
$ws1.Range("O10").Value = 1
$ws1.Range("P10").Value = 82.991389424999994
$ws1.Range("Q10").Value = 5944.0547263999997
$ws1.Range("R10").Value = 14
$ws1.Range("S10").Value = 69590.731713000001
$ws1.Range("T10").Value = 1479.9353234
$ws1.Range("U10").Value = 3
$ws1.Range("V10").Value = 12219.597841999999
$ws1.Range("W10").Value = 7423.9900497999997
$ws1.Range("X10").Value = 22
$ws1.Range("Y10").Value = 71764.104538
$ws1.Range("Z10").Value = 4231.8666855000001
$ws1.Range("AA10").Value = 16.7
$ws1.Range("AB10").Value = 32838.154939
$ws1.Range("AC10").Value = 5946.4651740999998
$ws1.Range("AD10").Value = 14
$ws1.Range("AE10").Value = 69590.661124999999
$ws1.Range("AF10").Value = 1482.3457711000001
$ws1.Range("AG10").Value = 4
$ws1.Range("AH10").Value = 12219.431757
$ws1.Range("AI10").Value = 7428.8109452999997
$ws1.Range("AJ10").Value = 22
$ws1.Range("AK10").Value = 71763.911078000005
$ws1.Range("AL10").Value = 4233.1761532999999
$ws1.Range("AM10").Value = 18.125
$ws1.Range("AN10").Value = 32838.041462000001
$ws1.Range("AO10").Value = 166.79104477999999
$ws1.Range("AP10").Value = 4
$ws1.Range("AQ10").Value = 1350.506079
$ws1.Range("AR10").Value = 147.37474101000001
$ws1.Range("AS10").Value = 3.6
$ws1.Range("AT10").Value = 1344.4942458
$ws1.Range("AU10").Value = 168.67910448000001
$ws1.Range("AV10").Value = 5
$ws1.Range("AW10").Value = 1350.4747789999999
$ws1.Range("AX10").Value = 147.82556536000001
$ws1.Range("AY10").Value = 4
$ws1.Range("AZ10").Value = 1344.4606898
$ws1.Range("BA10").Value = 1.0124378109000001
$ws1.Range("BB10").Value = 1
$ws1.Range("BC10").Value = 0.110967306
$ws1.Range("BD10").Value = 1.0074626866
$ws1.Range("BE10").Value = 1
$ws1.Range("BF10").Value = 0.086171144699999994
$ws1.Range("BG10").Value = 13.324094804
$ws1.Range("BH10").Value = 0
$ws1.Range("BI10").Value = 79.288798958000001
$ws1.Range("BJ10").Value = 19.814127280000001
$ws1.Range("BK10").Value = 0
$ws1.Range("BL10").Value = 159.85899377000001

# Sheet 2: "only branch segs gt 1" - add row 10 for godot.
# (D10, F10, G10 are formulas in this sheet - set below, after their
# precedent cells B10/C10/E10 already have values.)
$ws2.Range("A10").Value = "godot"
$ws2.Range("B10").Value = 402
$ws2.Range("C10").Value = 613
$ws2.Range("E10").Value = 311
$ws2.Range("H10").Value = 3.11
$ws2.Range("I10").Value = 3
$ws2.Range("J10").Value = 1.7459621999999999
$ws2.Range("K10").Value = 68.400000000000006
$ws2.Range("L10").Value = 5.5
$ws2.Range("M10").Value = 340.01175999999998
$ws2.Range("N10").Value = 20.143689999999999
$ws2.Range("O10").Value = 1.6666700000000001
$ws2.Range("P10").Value = 84.762940999999998
$ws2.Range("Q10").Value = 68.430000000000007
$ws2.Range("R10").Value = 5.5
$ws2.Range("S10").Value = 340.00637
$ws2.Range("T10").Value = 20.15569
$ws2.Range("U10").Value = 1.6667000000000001
$ws2.Range("V10").Value = 84.760571999999996
$ws2.Range("W10").Value = 14432.54
$ws2.Range("X10").Value = 172.5
$ws2.Range("Y10").Value = 129590.84
$ws2.Range("Z10").Value = 2002.15
$ws2.Range("AA10").Value = 25.5
$ws2.Range("AB10").Value = 11109.769
$ws2.Range("AC10").Value = 16434.689999999999
$ws2.Range("AD10").Value = 253
$ws2.Range("AE10").Value = 129944.67
$ws2.Range("AF10").Value = 3602.3541
$ws2.Range("AG10").Value = 69.5
$ws2.Range("AH10").Value = 22397.338
$ws2.Range("AI10").Value = 14442.23
$ws2.Range("AJ10").Value = 175
$ws2.Range("AK10").Value = 129590.05
$ws2.Range("AL10").Value = 2011.84
$ws2.Range("AM10").Value = 34.5
$ws2.Range("AN10").Value = 11108.566000000001
$ws2.Range("AO10").Value = 16454.07
$ws2.Range("AP10").Value = 266
$ws2.Range("AQ10").Value = 129942.88
$ws2.Range("AR10").Value = 3607.6181000000001
$ws2.Range("AS10").Value = 73.75
$ws2.Range("AT10").Value = 22396.812999999998
$ws2.Range("AU10").Value = 114.42
$ws2.Range("AV10").Value = 22.5
$ws2.Range("AW10").Value = 289.44914999999997
$ws2.Range("AX10").Value = 36.366458999999999
$ws2.Range("AY10").Value = 7.75
$ws2.Range("AZ10").Value = 72.621058000000005
$ws2.Range("BA10").Value = 122.01
$ws2.Range("BB10").Value = 31.5
$ws2.Range("BC10").Value = 290.16838000000001
$ws2.Range("BD10").Value = 38.178773
$ws2.Range("BE10").Value = 9.6666699999999999
$ws2.Range("BF10").Value = 72.885351999999997
$ws2.Range("BG10").Value = 1.05
$ws2.Range("BH10").Value = 1
$ws2.Range("BI10").Value = 0.21904290000000001
$ws2.Range("BJ10").Value = 1.03
$ws2.Range("BK10").Value = 1
$ws2.Range("BL10").Value = 0.1714466
$ws2.Range("BM10").Value = 53.562860999999998
$ws2.Range("BN10").Value = 3.1804199999999998
$ws2.Range("BO10").Value = 152.60174000000001
$ws2.Range("BP10").Value = 79.652792000000005
$ws2.Range("BQ10").Value = 4.83528
$ws2.Range("BR10").Value = 314.15884999999997
$ws2.Range("D10").Formula = "=B10-301"
$ws2.Range("F10").Formula = "=D10/B10"
$ws2.Range("G10").Formula = "=E10/C10"

# Update the saved selection on each sheet (sheet2 stays the active/tabSelected
# sheet, so select it last).
[void]$ws1.Range("C15").Select()
[void]$ws2.Range("E16").Select()
